# Update the "想去人数" (F) and "最低票价" (G) figures on the "展览" and
# "全部类型" sheets, which contain duplicated data (gh-pages output refresh).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 6670

    $ws.Range("F3").Value = 190
    $ws.Range("G3").Value = "不可售"

    $ws.Range("F6").Value = 2070

    $ws.Range("F7").Value = 1579

    $ws.Range("F9").Value = 1021

    $ws.Range("F10").Value = 460

    $ws.Range("F11").Value = 21

    $ws.Range("F12").Value = 5649
}
